$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the part info for the 8.2uH inductor (L4, row 17) with the
# SnapEDA-matched component (new Description / Digi-Key PN / MPN).
$ws.Range("C17").Value = "FIXED IND 8.2UH 1.3A 153.6 MOHM"
$ws.Range("D17").Value = "490-16115-1-ND"
$ws.Range("E17").Value = "LQH43PH8R2M26L"

# Reflect the user's last selection on the sheet.
$ws.Range("E17").Select()
